# Fill in the bill-form header fields (name, designation, year, term,
# department/branch, and the "in words" amount) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# নাম: (Name) -> add the teacher's name
$ws.Range("A3").Value = "নাম: Dr. S. M. Rabiul Alam (Hum)"

# পদবী: (Designation) -> add "অধ্যাপক" (Professor)
$ws.Range("A4").Value = "পদবী: অধ্যাপক"

# বর্ষ : (Year) value cell -> ৪র্থ (4th)
$ws.Range("G4").Value = "৪র্থ"

# টার্ম : (Term) value cell -> ১ম (1st)
$ws.Range("I4").Value = "১ম"

# বিভাগ/শাখা: (Department/Branch) value cell -> সিএসই (CSE)
$ws.Range("B5").Value = "সিএসই"

# বিভাগ : (Department) -> add "হুম" (Hum)
$ws.Range("F5").Value = "বিভাগ :হুম"

# কথায়: (In words) -> amount in words matching the 2700 total in I32
$ws.Range("A32").Value = "কথায়:দুই হাজার সাতশো টাকা মাত্র।"
